# 00_todo.xlsx - "2024-11" sheet (ActiveSheet): fill in the daily todo rows
# 15-26 (2024-11-09 .. 2024-11-20) with the standard recurring plan, and
# extend the remark for the last two days (19th/20th) to mention the
# YouTube policy distraction. Also move the viewport/selection down to
# where the new entries are (row 25 / D31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text blocks (shared across the recurring rows), matching existing
# shared-string content used elsewhere on this sheet.
$todoText    = "1、6：00 get up`n2、23：00 sleep`n3、workout"
$englishText = "1、300 words`n2、oral practise `n3、1 listen test`n4、0.5h attentive listen`n5、1h tech listen`n"
$techText    = "1、algorithm `n2、`n"
$remarkText  = "1、stare at stocks"
$remarkText2 = "1、stare at stocks`n2、waste in YouTube policy"

for ($r = 15; $r -le 24; $r++) {
    $ws.Cells.Item($r, 2).Value = $todoText
    $ws.Cells.Item($r, 2).Font.Strikethrough = $true

    $ws.Cells.Item($r, 3).Value = $englishText
    $ws.Cells.Item($r, 3).Font.Strikethrough = $true

    $ws.Cells.Item($r, 4).Value = $techText

    $ws.Cells.Item($r, 5).Value = $remarkText
}

for ($r = 25; $r -le 26; $r++) {
    $ws.Cells.Item($r, 2).Value = $todoText
    $ws.Cells.Item($r, 2).Font.Strikethrough = $true

    $ws.Cells.Item($r, 3).Value = $englishText
    $ws.Cells.Item($r, 3).Font.Strikethrough = $true

    $ws.Cells.Item($r, 4).Value = $techText

    $ws.Cells.Item($r, 5).Value = $remarkText2
}

# Row height for the newly-filled rows (wrapped multi-line content).
$ws.Range("A15:A26").EntireRow.RowHeight = 84

# Move selection / viewport to the newly-edited area.
$ws.Range("D31").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
